$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before "总计", matching the look of
#    the existing per-quarter holding sheets (e.g. "2021-Q4"), and fill it
#    in with the new quarter's fund-holding data.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: adding a sheet "before" the "总计" sheet re-seats the handle used to
# reference it, so re-resolve "总计" by name afterwards - don't reuse
# $totalSheetBefore for writes.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy formatting only (fonts/borders/alignment) from the 2021-Q4 sheet so
# the header row and index column match the sheets around it. Column A of
# the header row is untouched/blank in these sheets, so copy B:H and A:H
# separately to avoid manufacturing a spurious A1 cell.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$q4.Range("A2:H2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "070031"
$newSheet.Range("C2").Value = "嘉实全球房地产(QDII)"
$newSheet.Range("D2").Value = "0.60"
$newSheet.Range("E2").Value = "95.08"
$newSheet.Range("F2").Value = "2.65"
$newSheet.Range("G2").Value = "0.0159"
$newSheet.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new top data row for 2022-Q1 and
#    push the existing quarters down, renumbering the index column.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The freshly inserted row inherits stray formatting from its neighbours -
# clear it and restore just the index-column style (copied from the row
# below, which still carries it).
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.02

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
